$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 228
$ws.Range("F3").Value = 1414
$ws.Range("F4").Value = 19858
$ws.Range("F6").Value = 308
$ws.Range("F7").Value = 1094
$ws.Range("F8").Value = 13
$ws.Range("F9").Value = 7505
$ws.Range("F10").Value = 0
$ws.Range("F12").Value = 259
$ws.Range("F13").Value = 37
$ws.Range("F15").Value = 0
$ws.Range("F16").Value = 8
$ws.Range("F17").Value = 233
$ws.Range("F19").Value = 0
$ws.Range("F22").Value = 0
$ws.Range("F23").Value = 48
$ws.Range("F24").Value = 61
$ws.Range("F25").Value = 67
$ws.Range("F27").Value = 1093
$ws.Range("F29").Value = 0
$ws.Range("F31").Value = 5220
$ws.Range("F35").Value = 25
$ws.Range("F36").Value = 86
$ws.Range("F37").Value = 51
$ws.Range("F38").Value = 12586
$ws.Range("F39").Value = 1330
$ws.Range("F40").Value = 71
$ws.Range("F41").Value = 21
$ws.Range("F44").Value = 355
$ws.Range("F46").Value = 318

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F3").Value = 33

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 228
$ws.Range("F3").Value = 1414
$ws.Range("F4").Value = 0
$ws.Range("F5").Value = 794
$ws.Range("F6").Value = 308
$ws.Range("F7").Value = 1094
$ws.Range("F8").Value = 13
$ws.Range("F9").Value = 0
$ws.Range("F11").Value = 729
$ws.Range("F14").Value = 153
$ws.Range("F15").Value = 114
$ws.Range("F16").Value = 8
$ws.Range("F17").Value = 233
$ws.Range("F19").Value = 1334
$ws.Range("F20").Value = 400
$ws.Range("F23").Value = 48
$ws.Range("F24").Value = 61
$ws.Range("F26").Value = 318
$ws.Range("F27").Value = 1093
$ws.Range("F28").Value = 0
$ws.Range("F30").Value = 178
$ws.Range("F32").Value = 558
$ws.Range("F35").Value = 0
$ws.Range("F36").Value = 0
$ws.Range("F38").Value = 86
$ws.Range("F39").Value = 51
$ws.Range("F40").Value = 12586
$ws.Range("F41").Value = 1330
$ws.Range("F42").Value = 71
$ws.Range("F43").Value = 21
$ws.Range("F46").Value = 355
$ws.Range("F47").Value = 3985
$ws.Range("F48").Value = 0
$ws.Range("F49").Value = 93

Write-Host "Applied all F-column updates."
